$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range('D2')
$c.NumberFormat = "@"
$c.Value = '51.373.18'
$c.Style = "Normal"
$ws.Range('E2').Value = '  -0.83%  '
$c = $ws.Range('D3')
$c.NumberFormat = "@"
$c.Value = '2.775.03'
$c.Style = "Normal"
$ws.Range('E3').Value = '  -0.13%  '
$ws.Range('E4').Value = '  -0.01%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '352.81'
$c.Style = "Normal"
$ws.Range('E5').Value = '  -0.80%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '108.13'
$c.Style = "Normal"
$ws.Range('E6').Value = '  -0.98%  '
$ws.Range('E7').Value = '  -2.18%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range('E8').Value = '  -0.02%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.587'
$c.Style = "Normal"
$ws.Range('E9').Value = '  -1.21%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '39.62'
$c.Style = "Normal"
$ws.Range('E10').Value = '  -0.79%  '
$ws.Range('E11').Value = '  +3.00%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '0.0835'
$c.Style = "Normal"
$ws.Range('E12').Value = '  -1.74%  '
$c = $ws.Range('D13')
$c.NumberFormat = "@"
$c.Value = '20.05'
$c.Style = "Normal"
$ws.Range('E13').Value = '  +3.39%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '7.56'
$c.Style = "Normal"
$ws.Range('E14').Value = '  -0.54%  '
$c = $ws.Range('D15')
$c.NumberFormat = "@"
$c.Value = '3.209.55'
$c.Style = "Normal"
$ws.Range('E15').Value = '  +0.01%  '
$c = $ws.Range('D16')
$c.NumberFormat = "@"
$c.Value = '2.780.69'
$c.Style = "Normal"
$ws.Range('E16').Value = '  -0.43%  '
$c = $ws.Range('D17')
$c.NumberFormat = "@"
$c.Value = '0.922'
$c.Style = "Normal"
$ws.Range('E17').Value = '  -1.42%  '
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '51.337.13'
$c.Style = "Normal"
$ws.Range('E18').Value = '  -0.75%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '7.60'
$c.Style = "Normal"
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '3.10'
$c.Style = "Normal"
$ws.Range('E20').Value = '  -1.81%  '
$ws.Range('E22').Value = '  -1.47%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '69.84'
$c.Style = "Normal"
$ws.Range('E23').Value = '  -0.13%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '265.44'
$c.Style = "Normal"
$ws.Range('E24').Value = '  -3.25%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '2.71'
$c.Style = "Normal"
$ws.Range('E25').Value = '  -0.66%  '
$ws.Range('E26').Value = '  -0.06%  '
$c = $ws.Range('D27')
$c.NumberFormat = "@"
$c.Value = '26.02'
$c.Style = "Normal"
$ws.Range('E27').Value = '  -2.02%  '
$ws.Range('E28').Value = '  +11.68%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '10.22'
$c.Style = "Normal"
$ws.Range('E29').Value = '  +0.83%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '36.35'
$c.Style = "Normal"
$ws.Range('E30').Value = '  +6.94%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '6.17'
$c.Style = "Normal"
$ws.Range('E32').Value = '  +8.45%  '
$ws.Range('E33').Value = '  +0.54%  '
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '0.0453'
$c.Style = "Normal"
$ws.Range('E34').Value = '  -1.93%  '
$ws.Range('E35').Value = '  +5.01%  '
$ws.Range('E36').Value = '  -2.21%  '
$ws.Range('E37').Value = '  -0.03%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '18.41'
$c.Style = "Normal"
$ws.Range('E38').Value = '  +2.10%  '
$c = $ws.Range('D39')
$c.NumberFormat = "@"
$c.Value = '3.15'
$c.Style = "Normal"
$ws.Range('E39').Value = '  -2.24%  '
$ws.Range('E40').Value = '  -1.54%  '
$ws.Range('E41').Value = '  +0.99%  '
$ws.Range('E42').Value = '  -0.74%  '
$c = $ws.Range('D43')
$c.NumberFormat = "@"
$c.Value = '120.40'
$c.Style = "Normal"
$ws.Range('E43').Value = '  -0.86%  '
$c = $ws.Range('D44')
$c.NumberFormat = "@"
$c.Value = '22.05'
$c.Style = "Normal"
$ws.Range('E45').Value = '  -2.07%  '
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '2.100.96'
$c.Style = "Normal"
$ws.Range('E46').Value = '  +1.71%  '
$ws.Range('E47').Value = '  +0.72%  '
$ws.Range('E48').Value = '  +4.85%  '
$ws.Range('B49').Value = 'SEI'
$ws.Range('C49').Value = 'https://coinranking.com/coin/8nxCqs-uj+sei-sei'
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '0.905'
$c.Style = "Normal"
$ws.Range('E49').Value = '  -2.79%  '
$ws.Range('B50').Value = 'THORChain'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '5.41'
$c.Style = "Normal"
$ws.Range('E50').Value = '  -4.91%  '
$ws.Range('E51').Value = '  +8.53%  '
